# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Row -> (old value, new value) for sheet "展览" (sheet1.xml)
$changesExhibition = @{
    2  = 1332
    4  = 14624
    5  = 17767
    9  = 218
    15 = 45
    16 = 49
    21 = 77
    22 = 65
    24 = 7329
    26 = 4
    30 = 5870
    31 = 69
    35 = 228
    36 = 5107
}

# Row -> new value for sheet "全部类型" (sheet4.xml) - row numbers are offset
# by +1 starting at row 23 relative to "展览" because of an extra row.
$changesAllTypes = @{
    2  = 1332
    4  = 14624
    5  = 17767
    9  = 218
    15 = 45
    16 = 49
    21 = 77
    23 = 65
    25 = 7329
    27 = 4
    32 = 5870
    33 = 69
    37 = 228
    38 = 5107
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $changesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $changesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $changesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $changesAllTypes[$row]
}
